$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sehir")

# Row 5 used to hold a leftover "ankara" test row; correct it to the real
# exported city (tekirdağ, id 59).
$ws.Cells.Item(5, 1).Value2 = 59
$ws.Cells.Item(5, 2).Value2 = "tekirdağ"

# Rows 6 and 7 were stray test rows ("edirne" / "adana") left behind by the
# broken Export button; remove them now that it's fixed.
$ws.Rows.Item(7).Delete()
$ws.Rows.Item(6).Delete()
